# ---------------------------------------------------------------------------
# Applies the commit "Add files via upload" to LINDAO ZUÑIGA BRYAN JOSE.xlsx:
#   1. Zero-out column G (PRESUPUESTO) in "VENTA MENSUAL" for rows 2..56.
#   2. Add a new worksheet "CUMPLIMIENTO MENSUAL" (after "VENTA MENSUAL")
#      summarising budget vs. sales vs. compliance per GRUPO.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. VENTA MENSUAL -> zero the PRESUPUESTO column (G2:G56)
# ---------------------------------------------------------------------------
$ventaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$ventaMensual.Range("G2:G56").Value = 0

# ---------------------------------------------------------------------------
# 2. Add the new "CUMPLIMIENTO MENSUAL" sheet, placed after "VENTA MENSUAL"
#    (i.e. at the end of the tab strip).
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "CUMPLIMIENTO MENSUAL"

# Column widths (characters), matching the source workbook's other sheets.
$ws.Columns.Item(1).ColumnWidth = 25.166666666666668   # A -> 26
$ws.Columns.Item(2).ColumnWidth = 21.166666666666668   # B -> 22
$ws.Columns.Item(3).ColumnWidth = 21.166666666666668   # C -> 22
$ws.Columns.Item(4).ColumnWidth = 11.166666666666666   # D -> 12
$ws.Columns.Item(5).ColumnWidth = 21.166666666666668   # E -> 22
$ws.Columns.Item(6).ColumnWidth = 25.166666666666668   # F -> 26

# --- Header row (row 1): reuse the bold/bordered header style already used
#     by the other sheets by copying its formatting onto the new header row.
$ventasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ventasPorGrupo.Range("A1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1").Value = "ASESOR"
$ws.Range("B1").Value = "GRUPO"
$ws.Range("C1").Value = "PRESUPUESTO"
$ws.Range("D1").Value = "VENTA"
$ws.Range("E1").Value = "POR CUMPLIR"
$ws.Range("F1").Value = "CUMPLIMIENTO"

# --- Data rows (2..18): one per GRUPO, plus TOTAL on row 19.
$asesor = "LINDAO ZUÑIGA BRYAN JOSE"

$rows = @(
    @("240X120 PORCELANATO",  672.340305337043,   0,      672.340305337043,   0),
    @("240X80 PORCELANATO",   4992.1832,          0,      4992.1832,          0),
    @("FREGADEROS DE COCINA", 142.502095025027,   0,      142.502095025027,   0),
    @("GRANITO",              238.32,             0,      238.32,             0),
    @("GRIFERIAS",            106.82,             0,      106.82,             0),
    @("INODOROS",             2100,               0,      2100,               0),
    @("LAVABOS",              750,                0,      750,                0),
    @("LED",                  300,                0,      300,                0),
    @("NO RESURTIBLES",       650.25,             9.58,   640.67,             0.01473279507881584),
    @("OTROS",                0,                  0,      0,                  0),
    @("PANELES DECORATIVOS",  350,                0,      350,                0),
    @("PANELES PU",           230,                0,      230,                0),
    @("PANELES PVC",          483,                0,      483,                0),
    @("PIEDRA SINTERIZADA",   1505.12,            0,      1505.12,            0),
    @("PORCELANATO",          38417.17,           233.38, 38183.79,           0.006074887869148092),
    @("PUERTAS DE SEGURIDAD", 342,                0,      342,                0),
    @("SAL SOLUBLE",          4130,               0,      4130,               0)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $asesor
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $r = $r + 1
}

# Currency format ("$"#,##0.00, numFmtId 164) for PRESUPUESTO / VENTA / POR CUMPLIR
$ws.Range("C2:E18").NumberFormat = '"$"#,##0.00'
# Percentage format (0.00%, numFmtId 10) for CUMPLIMIENTO
$ws.Range("F2:F18").NumberFormat = "0.00%"

# --- TOTAL row (19)
$ws.Range("B19").Value = "TOTAL"
$ws.Range("B19").HorizontalAlignment = -4152   # xlRight

$ws.Range("C19").Value = 55409.70560036207
$ws.Range("D19").Value = 242.96
$ws.Range("E19").Value = 55166.74560036208
$ws.Range("F19").Value = 0.004384791389297914

$ws.Range("C19:E19").NumberFormat = '"$"#,##0.00'
$ws.Range("F19").NumberFormat = "0.00%"

$ws.Range("A1").Select()
